$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7, shifting existing rows 7.. down by one.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the "enable_year" parameter
$ws.Cells.Item(7, 1).Value = "CHE"
$ws.Cells.Item(7, 2).Value = "conv_elec_hydroror"
$ws.Cells.Item(7, 3).Value = "enable_year"
$ws.Cells.Item(7, 4).Value = "configuration"
$ws.Cells.Item(7, 7).Value = 1990

# Update selection to match target state
$ws.Range("H7").Select()
